$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-01 Saturday", "2025-11-02 Sunday"),
    @("17×49=833", "87×62=5394"),
    @("11×96=1056", "17×63=1071"),
    @("62×32=1984", "37×52=1924"),
    @("94×41=3854", "24×81=1944"),
    @("67×61=4087", "58×63=3654"),
    @("87×28=2436", "67×89=5963"),
    @("83×11=913", "81×47=3807"),
    @("42×16=672", "92×42=3864"),
    @("67×18=1206", "71×47=3337"),
    @("69×24=1656", "52×68=3536"),
    @("83×63=5229", "96×83=7968"),
    @("24×41=984", "92×33=3036"),
    @("47×93=4371", "38×61=2318"),
    @("51×74=3774", "30×29=870"),
    @("28×50=1400", "74×88=6512"),
    @("54×53=2862", "51×91=4641"),
    @("58×46=2668", "89×98=8722"),
    @("34×14=476", "89×49=4361"),
    @("96×23=2208", "16×47=752"),
    @("96×64=6144", "68×46=3128"),
    @("95×60=5700", "30×57=1710"),
    @("78×25=1950", "44×83=3652"),
    @("92×22=2024", "74×49=3626"),
    @("36×90=3240", "99×73=7227"),
    @("27×80=2160", "31×57=1767")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
